$d = $word.ActiveDocument

# Update the subtitle line
$d.Content.Find.Execute(
    "12 Prompts to Master Outlook with context and examples", $true, $false, $false, $false, $false,
    $true, 1, $false, "This document provides 12 prompts with explanations, examples, and tips.", 2
)

# Insert the new "Further reading" paragraph right after the last "Tips" line
# (Prompt 12's), before the closing copyright paragraph. Do this BEFORE the
# repeated Purpose/Prompt/Example Output/Tips replacements below so the
# "Tips: Adapt tone, check facts." text used to locate the paragraph is
# still unambiguous / unmodified at this point.
$lastTip = $d.Paragraphs | Where-Object { $_.Range.Text.TrimEnd("`r`a") -eq "Tips: Adapt tone, check facts." } | Select-Object -Last 1
$insertPoint = $d.Range($lastTip.Range.End, $lastTip.Range.End)
$insertPoint.InsertBefore("Further reading: https://hbr.org, https://www.microsoft.com/en-us/worklab`r")

# Update the repeated Purpose / Prompt / Example Output / Tips lines for all 12 prompt blocks
$d.Content.Find.Execute(
    "Purpose: Why this prompt is valuable.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Purpose: Save time triaging and drafting replies.", 2
)

$d.Content.Find.Execute(
    "Prompt text: [Insert Copilot prompt]", $true, $false, $false, $false, $false,
    $true, 1, $false, "Prompt: Example Outlook Copilot command.", 2
)

$d.Content.Find.Execute(
    "Example Output: [Sample Copilot result]", $true, $false, $false, $false, $false,
    $true, 1, $false, "Example Output: A short Copilot-generated email.", 2
)

$d.Content.Find.Execute(
    "Tips: Adapt tone, check facts.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Tip: Adjust tone for executive vs. peer audiences.", 2
)
